$wb = $excel.ActiveWorkbook

# --- Grade_1 sheet: mark row 4 (Grade1_(3).png) as sent + assign a Message_Id,
#     then remove the trailing 3 scheduled-but-unused rows (105-107). ---
$ws1 = $wb.Worksheets.Item("Grade_1")
$ws1.Cells.Item(4, 3).Value = 1
$ws1.Cells.Item(4, 4).Value = 7488
$ws1.Rows("105:107").Delete()

# --- Grade_2 sheet: mark row 4 (Grade2_(3).png) as sent + assign a Message_Id,
#     then append 4 new scheduled rows continuing the date/image sequence. ---
$ws2 = $wb.Worksheets.Item("Grade_2")
$ws2.Cells.Item(4, 3).Value = 1
$ws2.Cells.Item(4, 4).Value = 3461

$ws2.Cells.Item(99, 1).Value = "16,03,2023"
$ws2.Cells.Item(99, 2).Value = "Grade2_(98).png"
$ws2.Cells.Item(99, 3).Value = 0

$ws2.Cells.Item(100, 1).Value = "17,03,2023"
$ws2.Cells.Item(100, 2).Value = "Grade2_(99).png"
$ws2.Cells.Item(100, 3).Value = 0

$ws2.Cells.Item(101, 1).Value = "18,03,2023"
$ws2.Cells.Item(101, 2).Value = "Grade2_(100).png"
$ws2.Cells.Item(101, 3).Value = 0

$ws2.Cells.Item(102, 1).Value = "19,03,2023"
$ws2.Cells.Item(102, 2).Value = "Grade2_(101).png"
$ws2.Cells.Item(102, 3).Value = 0

# --- Grade_3 sheet: mark row 4 (Grade3_(3).png) as sent + assign a Message_Id,
#     then append 5 new scheduled rows continuing the date/image sequence. ---
$ws3 = $wb.Worksheets.Item("Grade_3")
$ws3.Cells.Item(4, 3).Value = 1
$ws3.Cells.Item(4, 4).Value = 2494

$ws3.Cells.Item(97, 1).Value = "14,03,2023"
$ws3.Cells.Item(97, 2).Value = "Grade3_(96).png"
$ws3.Cells.Item(97, 3).Value = 0

$ws3.Cells.Item(98, 1).Value = "15,03,2023"
$ws3.Cells.Item(98, 2).Value = "Grade3_(97).png"
$ws3.Cells.Item(98, 3).Value = 0

$ws3.Cells.Item(99, 1).Value = "16,03,2023"
$ws3.Cells.Item(99, 2).Value = "Grade3_(98).png"
$ws3.Cells.Item(99, 3).Value = 0

$ws3.Cells.Item(100, 1).Value = "17,03,2023"
$ws3.Cells.Item(100, 2).Value = "Grade3_(99).png"
$ws3.Cells.Item(100, 3).Value = 0

$ws3.Cells.Item(101, 1).Value = "18,03,2023"
$ws3.Cells.Item(101, 2).Value = "Grade3_(100).png"
$ws3.Cells.Item(101, 3).Value = 0

# --- Grade_4 sheet: mark row 4 (Grade4_(3).png) as sent + assign a Message_Id,
#     then append 7 new scheduled rows continuing the date/image sequence. ---
$ws4 = $wb.Worksheets.Item("Grade_4")
$ws4.Cells.Item(4, 3).Value = 1
$ws4.Cells.Item(4, 4).Value = 6226

$ws4.Cells.Item(96, 1).Value = "13,03,2023"
$ws4.Cells.Item(96, 2).Value = "Grade4_(95).png"
$ws4.Cells.Item(96, 3).Value = 0

$ws4.Cells.Item(97, 1).Value = "14,03,2023"
$ws4.Cells.Item(97, 2).Value = "Grade4_(96).png"
$ws4.Cells.Item(97, 3).Value = 0

$ws4.Cells.Item(98, 1).Value = "15,03,2023"
$ws4.Cells.Item(98, 2).Value = "Grade4_(97).png"
$ws4.Cells.Item(98, 3).Value = 0

$ws4.Cells.Item(99, 1).Value = "16,03,2023"
$ws4.Cells.Item(99, 2).Value = "Grade4_(98).png"
$ws4.Cells.Item(99, 3).Value = 0

$ws4.Cells.Item(100, 1).Value = "17,03,2023"
$ws4.Cells.Item(100, 2).Value = "Grade4_(99).png"
$ws4.Cells.Item(100, 3).Value = 0

$ws4.Cells.Item(101, 1).Value = "18,03,2023"
$ws4.Cells.Item(101, 2).Value = "Grade4_(100).png"
$ws4.Cells.Item(101, 3).Value = 0

$ws4.Cells.Item(102, 1).Value = "19,03,2023"
$ws4.Cells.Item(102, 2).Value = "Grade4_(101).png"
$ws4.Cells.Item(102, 3).Value = 0
